$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row of session data (Session 11, second strategy)
$ws.Range("A6").Value = 11
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = 4
$ws.Range("D6").Value = 2
$ws.Range("E6").Value = 4

# Move the active selection to F6 as in the updated workbook
$ws.Range("F6").Select()
